$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B44: change from text "3" to numeric 3
$ws.Range("B44").Value = 3

# Add new row 45 with data (B45 must remain text "3", not numeric)
$ws.Range("A45").Value = "Ruilin"
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "3"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "无"
$ws.Range("D45").Value = "DIS"
$ws.Range("E45").Value = "RES"
$ws.Range("F45").Value = "18e2478f-5f8b-460a-bbaf-4b86b95999fd"
$ws.Range("G45").Value = "B1IDRdeCW_annotated.xlsx"
$ws.Range("H45").Value = "This paper presents three observations to understand binary network in Courbariaux, Hubara et al. (2016)."
